$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("ECs", "Ccl21b", "Ccr10", "ECs", 1, 0.3333333333333333, 0.004739, 0.014217, 0.02588570741885795, 0.02588570741885795, 3, 1, 1.058783666666667, 3.176351, 0.5014862149947701, 0.5014862149947702, 0.005017575796333334, 0.04515818216700001, 0.01298132543594511, 0.01298132543594512),
  @("ECs", "Ccl21b", "Ccr10", "FAPs", 1, 0.3333333333333333, 0.004739, 0.014217, 0.02588570741885795, 0.02588570741885795, 2, 0.6666666666666666, 1.003554, 3.010662, 0.475327031240749, 0.4753270312407492, 0.004755842406000001, 0.042802581654, 0.01230417645897238, 0.01230417645897238),
  @("ECs", "Ccl21b", "Ccr10", "Resolving-Mac", 1, 0.3333333333333333, 0.004739, 0.014217, 0.02588570741885795, 0.02588570741885795, 1, 0.3333333333333333, 0.048954, 0.146862, 0.02318675376448066, 0.02318675376448067, 0.000231993006, 0.002087937054, 0.0006002055239404495, 0.0006002055239404497),
  @("FAPs", "Ccl21b", "Ccr10", "ECs", 2, 0.6666666666666666, 0.134289, 0.402867, 0.733523056250478, 0.733523056250478, 3, 1, 1.058783666666667, 3.176351, 0.5014862149947701, 0.5014862149947702, 0.142182999813, 1.279646998317, 0.367851701090448, 0.3678517010904481),
  @("FAPs", "Ccl21b", "Ccr10", "FAPs", 2, 0.6666666666666666, 0.134289, 0.402867, 0.733523056250478, 0.733523056250478, 2, 0.6666666666666666, 1.003554, 3.010662, 0.475327031240749, 0.4753270312407492, 0.134766263106, 1.212896367954, 0.3486633366741806, 0.3486633366741808),
  @("FAPs", "Ccl21b", "Ccr10", "Resolving-Mac", 2, 0.6666666666666666, 0.134289, 0.402867, 0.733523056250478, 0.733523056250478, 1, 0.3333333333333333, 0.048954, 0.146862, 0.02318675376448066, 0.02318675376448067, 0.006573983706000001, 0.059165853354, 0.01700801848584913, 0.01700801848584913),
  @("MuSCs", "Ccl21b", "Ccr10", "ECs", 2, 0.6666666666666666, 0.04404599999999999, 0.132138, 0.2405912363306641, 0.2405912363306641, 3, 1, 1.058783666666667, 3.176351, 0.5014862149947701, 0.5014862149947702, 0.046635185382, 0.419716668438, 0.1206531884683769, 0.120653188468377),
  @("MuSCs", "Ccl21b", "Ccr10", "FAPs", 2, 0.6666666666666666, 0.04404599999999999, 0.132138, 0.2405912363306641, 0.2405912363306641, 2, 0.6666666666666666, 1.003554, 3.010662, 0.475327031240749, 0.4753270312407492, 0.044202539484, 0.3978228553559999, 0.114359518107596, 0.114359518107596),
  @("MuSCs", "Ccl21b", "Ccr10", "Resolving-Mac", 2, 0.6666666666666666, 0.04404599999999999, 0.132138, 0.2405912363306641, 0.2405912363306641, 1, 0.3333333333333333, 0.048954, 0.146862, 0.02318675376448066, 0.02318675376448067, 0.002156227884, 0.019406050956, 0.005578529754691082, 0.005578529754691083)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $i + 2
    for ($j = 0; $j -lt $row.Length; $j++) {
        $c = $j + 1
        $ws.Cells.Item($r, $c).Value = $row[$j]
    }
}

Write-Output "done"